$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.949.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.293.63"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.82%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.316.08"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.12"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.703.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.955.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("E17").Value = "  +1.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.308.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.31%  "

$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("E22").Value = "  +4.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("E29").Value = "  +3.16%  "

$ws.Range("E30").Value = "  +2.28%  "

$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("E32").Value = "  +4.31%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.36%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("E36").Value = "  -4.19%  "

$ws.Range("E37").Value = "  +2.61%  "

$ws.Range("E38").Value = "  +3.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.83"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.377"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "133.96"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.94%  "

$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("E44").Value = "  -2.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "261.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0505"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0914"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.75%  "

$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.72%  "
